# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
# This script re-orders several match rows on the "Northern Ireland Premier"
# sheet. For each group of rows listed below, every row's data (columns
# B, and F through AC) is rotated "upwards": the new content of the first
# row in the group becomes the old content of the second row, and so on,
# with the last row's old content wrapping around to the first row.
# Columns A (id), C, D and E (Div, Div Original Name, Date) are left
# untouched, since they are identical across each group of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($ws, $row) {
    $b = $ws.Range("B$row").Value2
    $fac = $ws.Range("F$row`:AC$row").Value2
    return @{ B = $b; FAC = $fac }
}

function Set-RowData($ws, $row, $data) {
    $ws.Range("B$row").Value2 = $data.B
    $ws.Range("F$row`:AC$row").Value2 = $data.FAC
}

function Rotate-Rows($ws, [int[]]$rows) {
    # Capture the current data for every row in the group first, so that
    # later writes don't clobber data we still need to read.
    $data = @{}
    foreach ($r in $rows) {
        $data[$r] = Get-RowData $ws $r
    }

    # new(row[i]) = old(row[i+1]), wrapping around at the end.
    $count = $rows.Count
    for ($i = 0; $i -lt $count; $i++) {
        $targetRow = $rows[$i]
        $sourceRow = $rows[($i + 1) % $count]
        Set-RowData $ws $targetRow $data[$sourceRow]
    }
}

Rotate-Rows $ws @(105, 106, 107)
Rotate-Rows $ws @(121, 122)
Rotate-Rows $ws @(162, 163)
Rotate-Rows $ws @(170, 171)
Rotate-Rows $ws @(190, 191, 192)
